$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.601.22'
$ws.Range('E2').Value = '  -0.63%  '
$ws.Range('D3').Value = '3.899.69'
$ws.Range('E3').Value = '  +2.39%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'601.85"
$ws.Range('E5').Value = '  +0.08%  '
$ws.Range('D6').Value = "'166.97"
$ws.Range('E6').Value = '  +1.96%  '
$ws.Range('D7').Value = '3.898.56'
$ws.Range('E7').Value = '  +2.38%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('E9').Value = '  -1.16%  '
$ws.Range('D10').Value = "'0.169"
$ws.Range('E10').Value = '  -0.72%  '
$ws.Range('D11').Value = "'6.46"
$ws.Range('E11').Value = '  +2.33%  '
$ws.Range('D12').Value = "'0.462"
$ws.Range('E12').Value = '  +0.25%  '
$ws.Range('D13').Value = "'0.0000256"
$ws.Range('E13').Value = '  +4.06%  '
$ws.Range('D14').Value = "'37.35"
$ws.Range('E14').Value = '  +0.39%  '
$ws.Range('D15').Value = '4.551.06'
$ws.Range('E15').Value = '  +2.30%  '
$ws.Range('D16').Value = '3.895.36'
$ws.Range('E16').Value = '  +2.66%  '
$ws.Range('D17').Value = '68.647.60'
$ws.Range('E17').Value = '  -0.75%  '
$ws.Range('D18').Value = "'7.48"
$ws.Range('E18').Value = '  +0.10%  '
$ws.Range('D19').Value = "'17.40"
$ws.Range('E19').Value = '  +0.85%  '
$ws.Range('E20').Value = '  -2.27%  '
$ws.Range('E21').Value = '  -4.19%  '
$ws.Range('D22').Value = "'491.06"
$ws.Range('E22').Value = '  +1.11%  '
$ws.Range('D23').Value = "'0.727"
$ws.Range('E23').Value = '  +0.86%  '
$ws.Range('D24').Value = "'0.0000165"
$ws.Range('E24').Value = '  +2.95%  '
$ws.Range('D25').Value = "'84.83"
$ws.Range('E25').Value = '  +0.08%  '
$ws.Range('D26').Value = "'2.22"
$ws.Range('E26').Value = '  -1.36%  '
$ws.Range('D27').Value = "'12.04"
$ws.Range('E27').Value = '  -1.73%  '
$ws.Range('D28').Value = "'10.15"
$ws.Range('E28').Value = '  +1.20%  '
$ws.Range('E29').Value = '  +0.11%  '
$ws.Range('D30').Value = "'2.95"
$ws.Range('E30').Value = '  -0.66%  '
$ws.Range('D31').Value = '4.050.18'
$ws.Range('E31').Value = '  +2.37%  '
$ws.Range('E32').Value = '  -1.11%  '
$ws.Range('D33').Value = "'7.70"
$ws.Range('E33').Value = '  -4.09%  '
$ws.Range('D34').Value = "'31.79"
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('D35').Value = '3.856.78'
$ws.Range('E35').Value = '  +2.76%  '
$ws.Range('E36').Value = '  -0.32%  '
$ws.Range('E37').Value = '  +0.73%  '
$ws.Range('D38').Value = "'0.140"
$ws.Range('E38').Value = '  -0.75%  '
$ws.Range('D39').Value = "'5.94"
$ws.Range('E39').Value = '  +0.84%  '
$ws.Range('D40').Value = "'3.20"
$ws.Range('E40').Value = '  +5.97%  '
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('E42').Value = '  -1.04%  '
$ws.Range('D43').Value = "'430.22"
$ws.Range('E43').Value = '  +0.30%  '
$ws.Range('E44').Value = '  -0.42%  '
$ws.Range('D45').Value = "'48.03"
$ws.Range('E45').Value = '  -1.19%  '
$ws.Range('E46').Value = '  +2.19%  '
$ws.Range('E47').Value = '  +0.01%  '
$ws.Range('D48').Value = "'0.000276"
$ws.Range('E48').Value = '  +21.04%  '
$ws.Range('D49').Value = "'143.03"
$ws.Range('E49').Value = '  +0.84%  '
$ws.Range('D50').Value = '2.802.68'
$ws.Range('E50').Value = '  -0.93%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = "'25.79"
$ws.Range('E51').Value = '  +3.58%  '
